$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '43.403.94'
$ws.Range('E2').Value = '  +3.65%  '

$ws.Range('D3').Value = '2.310.04'
$ws.Range('E3').Value = '  +3.06%  '

$ws.Range('E4').Value = '  -0.10%  '

$ws.Range('B5').Value = 'Solana'
$ws.Range('C5').Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '105.31'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +10.00%  '

$ws.Range('B6').Value = 'BNB'
$ws.Range('C6').Value = 'https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb'
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '308.52'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +0.82%  '

$ws.Range('E7').Value = '  +0.79%  '

$ws.Range('E8').Value = '  -0.06%  '

$ws.Range('E9').Value = '  +6.22%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '35.89'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +4.25%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '52.64'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +2.60%  '

$ws.Range('E12').Value = '  -0.05%  '

$ws.Range('E13').Value = '  -1.04%  '

$ws.Range('E14').Value = '  +3.47%  '

$ws.Range('D15').Value = '2.669.64'
$ws.Range('E15').Value = '  +2.93%  '

$ws.Range('E16').Value = '  +5.62%  '

$ws.Range('D17').Value = '2.310.75'
$ws.Range('E17').Value = '  +2.64%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.800'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +3.09%  '

$ws.Range('D19').Value = '43.354.02'
$ws.Range('E19').Value = '  +3.67%  '

$ws.Range('D20').Value = '0.0₃0921'
$ws.Range('E20').Value = '  +2.90%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '11.84'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -1.98%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.19'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +5.04%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '67.94'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +1.39%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '240.28'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +2.31%  '

$ws.Range('E25').Value = '  +4.40%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.59'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +1.44%  '

$ws.Range('E27').Value = '  +0.19%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '24.98'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +8.17%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.25'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +6.64%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '36.23'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -3.61%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '9.59'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +1.83%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '162.32'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -2.91%  '

$ws.Range('E33').Value = '  +1.94%  '

$ws.Range('E34').Value = '  -0.10%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '18.23'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +5.30%  '

$ws.Range('E36').Value = '  +6.41%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.0732'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +2.41%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '4.61'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +15.14%  '

$ws.Range('E39').Value = '  -0.97%  '

$ws.Range('E40').Value = '  +3.74%  '

$ws.Range('E41').Value = '  +4.53%  '

$ws.Range('E42').Value = '  +0.68%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '2.50'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +15.30%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.0289'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +3.34%  '

$ws.Range('D45').Value = '1.962.99'
$ws.Range('E45').Value = '  +1.50%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '18.84'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +2.02%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '3.04'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +6.07%  '

$ws.Range('E48').Value = '  +7.07%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '57.72'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +7.94%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '2.96'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +3.14%  '

$ws.Range('E51').Value = '  +8.42%  '
